# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" (E16:E23) column used to list the debt periods in
# ascending order (2010, 2011, 2012, 2101, 2102, 2103, 2104, 2105).
# The database refresh replaces that listing with the same eight periods
# but in descending (most-recent-first) order, and the one differing
# "Valor Mora" amount (which used to sit on the last/oldest period) now
# sits on the first/newest period instead - i.e. the two distinct amounts
# in column F (28090 and 35112) swap rows along with the reordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("2105", "2104", "2103", "2102", "2101", "2012", "2011", "2010")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

$ws.Cells.Item(16, 6).Value = 28090
$ws.Cells.Item(23, 6).Value = 35112
